$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
  @("38+18=56", "90+6=96", "58-44=14", "3+12=15", "31+34=65"),
  @("30+69=99", "71-66=5", "69-49=20", "54-46=8", "89-50=39"),
  @("26+26=52", "23+65=88", "7+41=48", "57-52=5", "12+39=51"),
  @("1+73=74", "57-36=21", "49-0=49", "62+1=63", "76-69=7"),
  @("21+38=59", "92-27=65", "9-3=6", "9+85=94", "65-50=15"),
  @("6+30=36", "19+50=69", "35+53=88", "55+35=90", "11+78=89"),
  @("56+12=68", "52-26=26", "39+12=51", "6+34=40", "7+58=65"),
  @("97-53=44", "50-46=4", "66-59=7", "16-2=14", "11-10=1"),
  @("58+27=85", "92-29=63", "21+61=82", "3+42=45", "4+1=5"),
  @("36+23=59", "41+56=97", "23-20=3", "84+3=87", "18+37=55"),
  @("85-49=36", "52-3=49", "62-18=44", "66-46=20", "11+54=65"),
  @("60+15=75", "60+34=94", "9+86=95", "72-25=47", "44+17=61"),
  @("33-15=18", "2+60=62", "87-2=85", "17+75=92", "81-10=71"),
  @("3+94=97", "36+42=78", "11+22=33", "14+38=52", "93-62=31"),
  @("34-7=27", "31+26=57", "30-20=10", "11+53=64", "99-20=79"),
  @("44-25=19", "56-44=12", "61+9=70", "56+26=82", "66-60=6"),
  @("91-42=49", "70+24=94", "80+9=89", "15+11=26", "48-31=17"),
  @("91-62=29", "32+18=50", "71-13=58", "70+11=81", "71-66=5"),
  @("60-59=1", "43+5=48", "32-12=20", "95-74=21", "70+9=79"),
  @("79-39=40", "75-46=29", "57-18=39", "21+7=28", "87+3=90")
)
for ($r = 1; $r -le $values.Count; $r++) {
  $rowVals = $values[$r - 1]
  for ($c = 1; $c -le $rowVals.Count; $c++) {
    $t.Cell($r, $c).Range.Text = $rowVals[$c - 1]
  }
}
Write-Host "done"
